# Scheduled runner update: refresh Universalis market-price snapshots
# (currentAveragePrice* / LevePrice* / LeveProfit*) across all Job sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1700
$ws.Range("J43").Value = 2120
$ws.Range("L43").Value = 2120
$ws.Range("N43").Value = -2258
$ws.Range("H51").Value = 50002
$ws.Range("J51").Value = 50002
$ws.Range("L51").Value = 50002
$ws.Range("N51").Value = -50970
$ws.Range("H52").Value = 1200
$ws.Range("I52").Value = 400
$ws.Range("J52").Value = 2000
$ws.Range("K52").Value = 1200
$ws.Range("L52").Value = 6000
$ws.Range("M52").Value = -1040
$ws.Range("N52").Value = -6320
$ws.Range("H53").Value = 1891.1666
$ws.Range("I53").Value = 3105.3333
$ws.Range("K53").Value = 3105.3333
$ws.Range("M53").Value = -2468.3333
$ws.Range("H101").Value = 337.9091
$ws.Range("I101").Value = 303.42856
$ws.Range("J101").Value = 398.25
$ws.Range("K101").Value = 910.28568
$ws.Range("L101").Value = 1194.75
$ws.Range("M101").Value = 711.71432
$ws.Range("N101").Value = -4438.75
$ws.Range("H127").Value = 2400
$ws.Range("I127").Value = 2100
$ws.Range("K127").Value = 6300
$ws.Range("M127").Value = -1340
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("H138").Value = 3555.2222
$ws.Range("J138").Value = 3933
$ws.Range("L138").Value = 11799
$ws.Range("N138").Value = -22079

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2562.6667
$ws.Range("I45").Value = 2565
$ws.Range("K45").Value = 2565
$ws.Range("M45").Value = -2188
$ws.Range("H61").Value = 2965.3333
$ws.Range("I61").Value = 2498.5
$ws.Range("J61").Value = 3899
$ws.Range("K61").Value = 2498.5
$ws.Range("L61").Value = 3899
$ws.Range("M61").Value = -2286.5
$ws.Range("N61").Value = -4323
$ws.Range("H74").Value = 2249.75
$ws.Range("I74").Value = 2500
$ws.Range("J74").Value = 1999.5
$ws.Range("K74").Value = 2500
$ws.Range("L74").Value = 1999.5
$ws.Range("M74").Value = -1626
$ws.Range("N74").Value = -3747.5
$ws.Range("H77").Value = 2249.75
$ws.Range("I77").Value = 2500
$ws.Range("J77").Value = 1999.5
$ws.Range("K77").Value = 12500
$ws.Range("L77").Value = 9997.5
$ws.Range("M77").Value = -8132
$ws.Range("N77").Value = -18733.5
$ws.Range("H97").Value = 500
$ws.Range("I97").Value = 500
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 500
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -4
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 3400
$ws.Range("I122").Value = 3400
$ws.Range("K122").Value = 10200
$ws.Range("M122").Value = -7750
$ws.Range("H132").Value = 2475.3572
$ws.Range("I132").Value = 2475.3572
$ws.Range("K132").Value = 7426.071599999999
$ws.Range("M132").Value = -4896.071599999999
$ws.Range("H136").Value = 2965.3333
$ws.Range("I136").Value = 2498.5
$ws.Range("J136").Value = 3899
$ws.Range("K136").Value = 7495.5
$ws.Range("L136").Value = 11697
$ws.Range("M136").Value = -4945.5
$ws.Range("N136").Value = -16797

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3000.3333
$ws.Range("I4").Value = 2001
$ws.Range("J4").Value = 3500
$ws.Range("K4").Value = 2001
$ws.Range("L4").Value = 3500
$ws.Range("M4").Value = -1889
$ws.Range("N4").Value = -3724
$ws.Range("H7").Value = 194.29411
$ws.Range("I7").Value = 90.5
$ws.Range("J7").Value = 286.55554
$ws.Range("K7").Value = 90.5
$ws.Range("L7").Value = 286.55554
$ws.Range("M7").Value = 22.5
$ws.Range("N7").Value = -512.5555400000001
$ws.Range("H132").Value = 2864.2222
$ws.Range("I132").Value = 2864.2222
$ws.Range("K132").Value = 8592.6666
$ws.Range("M132").Value = -6062.6666

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 240.3125
$ws.Range("I2").Value = 118.333336
$ws.Range("J2").Value = 397.14285
$ws.Range("K2").Value = 710.000016
$ws.Range("L2").Value = 2382.8571
$ws.Range("M2").Value = -597.000016
$ws.Range("N2").Value = -2608.8571
$ws.Range("H29").Value = 136.5
$ws.Range("I29").Value = 84.57143000000001
$ws.Range("K29").Value = 253.71429
$ws.Range("M29").Value = 23.28570999999999
$ws.Range("H139").Value = 2155.4285
$ws.Range("I139").Value = 619
$ws.Range("K139").Value = 1857
$ws.Range("M139").Value = 3283

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 70834.164
$ws.Range("J5").Value = 1005
$ws.Range("L5").Value = 1005
$ws.Range("N5").Value = -1229
$ws.Range("H7").Value = 700600
$ws.Range("J7").Value = 501500
$ws.Range("L7").Value = 501500
$ws.Range("N7").Value = -501724
$ws.Range("H8").Value = 700600
$ws.Range("J8").Value = 501500
$ws.Range("L8").Value = 501500
$ws.Range("N8").Value = -501778
$ws.Range("H17").Value = 289.5
$ws.Range("I17").Value = 154
$ws.Range("J17").Value = 425
$ws.Range("K17").Value = 154
$ws.Range("L17").Value = 425
$ws.Range("M17").Value = 14
$ws.Range("N17").Value = -761
$ws.Range("H23").Value = 5368.1665
$ws.Range("I23").Value = 112
$ws.Range("J23").Value = 6419.4
$ws.Range("K23").Value = 112
$ws.Range("L23").Value = 6419.4
$ws.Range("M23").Value = 111
$ws.Range("N23").Value = -6865.4
$ws.Range("H32").Value = 35290
$ws.Range("J32").Value = 35290
$ws.Range("L32").Value = 35290
$ws.Range("N32").Value = -35882
$ws.Range("H122").Value = 2348.1428
$ws.Range("I122").Value = 1646.75
$ws.Range("K122").Value = 4940.25
$ws.Range("M122").Value = -2490.25
$ws.Range("H124").Value = 40000
$ws.Range("I124").Value = 40000
$ws.Range("K124").Value = 40000
$ws.Range("M124").Value = -35090

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6363
$ws.Range("I7").Value = 5999.3
$ws.Range("K7").Value = 5999.3
$ws.Range("M7").Value = -5887.3
$ws.Range("H126").Value = 6363
$ws.Range("I126").Value = 5999.3
$ws.Range("K126").Value = 17997.9
$ws.Range("M126").Value = -15527.9
$ws.Range("H136").Value = 2917.3794
$ws.Range("I136").Value = 2552.476
$ws.Range("K136").Value = 7657.428
$ws.Range("M136").Value = -5107.428

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 60000
$ws.Range("J27").Value = 60000
$ws.Range("L27").Value = 60000
$ws.Range("N27").Value = -60138
$ws.Range("H96").Value = 1006.6667
$ws.Range("I96").Value = 734.5454999999999
$ws.Range("J96").Value = 4000
$ws.Range("K96").Value = 734.5454999999999
$ws.Range("L96").Value = 4000
$ws.Range("M96").Value = 638.4545000000001
$ws.Range("N96").Value = -6746
$ws.Range("H115").Value = 60000
$ws.Range("J115").Value = 60000
$ws.Range("L115").Value = 60000
$ws.Range("N115").Value = -63134
$ws.Range("H122").Value = 3255.5557
$ws.Range("I122").Value = 2399
$ws.Range("J122").Value = 3683.8333
$ws.Range("K122").Value = 7197
$ws.Range("L122").Value = 11051.4999
$ws.Range("M122").Value = -4747
$ws.Range("N122").Value = -15951.4999
$ws.Range("H126").Value = 35795.438
$ws.Range("I126").Value = 31515.2
$ws.Range("K126").Value = 94545.60000000001
$ws.Range("M126").Value = -92075.60000000001
$ws.Range("H136").Value = 3777.76
$ws.Range("J136").Value = 3998
$ws.Range("L136").Value = 11994
$ws.Range("N136").Value = -17094

Write-Host "Scheduled runner: updated 202 price cells across 8 sheets."